# Applies updated crypto price/volume figures to columns D (Price) and E (Volume(1h))
# for rows 2-51, matching the refreshed data pulled by the scheduled GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '38.321.47'
$ws.Range('E2').Value = '  +1.13%  '
$ws.Range('D3').Value = '2.103.07'
$ws.Range('E3').Value = '  +3.22%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '''229.13'
$ws.Range('E5').Value = '  +0.33%  '
$ws.Range('D6').Value = '''0.613'
$ws.Range('E6').Value = '  +0.32%  '
$ws.Range('D7').Value = '''61.18'
$ws.Range('E7').Value = '  +1.01%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = '''0.381'
$ws.Range('E9').Value = '  +0.27%  '
$ws.Range('D10').Value = '''0.0852'
$ws.Range('E10').Value = '  +3.69%  '
$ws.Range('E11').Value = '  +0.09%  '
$ws.Range('D12').Value = '2.405.01'
$ws.Range('D13').Value = '''14.76'
$ws.Range('E13').Value = '  +1.62%  '
$ws.Range('D14').Value = '''22.40'
$ws.Range('E14').Value = '  +5.51%  '
$ws.Range('D15').Value = '''5.49'
$ws.Range('E15').Value = '  +6.01%  '
$ws.Range('D16').Value = '''0.779'
$ws.Range('E16').Value = '  +2.28%  '
$ws.Range('D17').Value = '2.093.70'
$ws.Range('E17').Value = '  +2.86%  '
$ws.Range('D18').Value = '38.240.84'
$ws.Range('E18').Value = '  +1.05%  '
$ws.Range('D19').Value = '''6.04'
$ws.Range('E19').Value = '  +2.25%  '
$ws.Range('D20').Value = '''70.41'
$ws.Range('E20').Value = '  +0.75%  '
$ws.Range('D21').Value = '0.0₃0836'
$ws.Range('E21').Value = '  +1.34%  '
$ws.Range('D22').Value = '''224.55'
$ws.Range('E22').Value = '  +0.24%  '
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').Value = '''2.44'
$ws.Range('E24').Value = '  +0.96%  '
$ws.Range('E25').Value = '  +2.36%  '
$ws.Range('D26').Value = '''169.99'
$ws.Range('E26').Value = '  +1.58%  '
$ws.Range('E27').Value = '  +0.78%  '
$ws.Range('D28').Value = '''0.131'
$ws.Range('E28').Value = '  +0.48%  '
$ws.Range('D29').Value = '''19.04'
$ws.Range('E29').Value = '  +0.84%  '
$ws.Range('E30').Value = '  +6.58%  '
$ws.Range('E31').Value = '  -0.54%  '
$ws.Range('E32').Value = '  +7.32%  '
$ws.Range('D33').Value = '''4.71'
$ws.Range('E33').Value = '  +4.19%  '
$ws.Range('D34').Value = '''4.44'
$ws.Range('E34').Value = '  +0.58%  '
$ws.Range('D35').Value = '''0.0608'
$ws.Range('E35').Value = '  +0.17%  '
$ws.Range('D36').Value = '''2.40'
$ws.Range('E36').Value = '  +5.03%  '
$ws.Range('D37').Value = '''6.44'
$ws.Range('E37').Value = '  +1.40%  '
$ws.Range('D38').Value = '''3.53'
$ws.Range('E38').Value = '  +5.87%  '
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('D40').Value = '''18.26'
$ws.Range('E40').Value = '  +3.02%  '
$ws.Range('D41').Value = '1.552.38'
$ws.Range('E41').Value = '  +0.74%  '
$ws.Range('D42').Value = '''100.21'
$ws.Range('E42').Value = '  +3.78%  '
$ws.Range('D43').Value = '''0.0220'
$ws.Range('E43').Value = '  +0.84%  '
$ws.Range('E44').Value = '  +0.74%  '
$ws.Range('D45').Value = '''0.0913'
$ws.Range('E45').Value = '  -0.08%  '
$ws.Range('D46').Value = '''4.16'
$ws.Range('E46').Value = '  +3.71%  '
$ws.Range('E47').Value = '  +1.24%  '
$ws.Range('D48').Value = '''7.45'
$ws.Range('E48').Value = '  +3.65%  '
$ws.Range('E49').Value = '  +1.64%  '
$ws.Range('D50').Value = '''3.00'
$ws.Range('E50').Value = '  +1.31%  '
$ws.Range('D51').Value = '2.296.02'
$ws.Range('E51').Value = '  +3.06%  '
